# The package carries three stray SharePoint/OneDrive "document library"
# CustomXML parts -- a content-type schema, a library form-template stub,
# and a library properties stub -- that the SharePoint sync client injects
# into every file it touches. They aren't referenced anywhere in the
# visible document content (no content control / XML mapping points at
# them), so as part of refreshing the Documentserver environment they are
# cleaned out of the package and the file is saved as a plain .docx again.

$d = $word.ActiveDocument

# The three stray parts, identified the way Word itself would distinguish
# them: by namespace URI (root element namespace of each part) and,
# failing that, by their storage GUID (the ds:itemID used in the matching
# itemProps#.xml).
$strayNamespaces = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms",
    "http://schemas.microsoft.com/office/2006/metadata/properties"
)
$strayIds = @(
    "{41635B2E-1A12-485A-9F81-7F501C6BE053}",
    "{0A532DDA-80AD-4477-AE1D-2EE09C89E47B}",
    "{CEDBCB29-50AD-4601-A4D2-369694091BD6}"
)

$parts = $d.CustomXMLParts

# Walk backwards so deleting an item doesn't shift the indices of the
# ones still to be examined.
for ($i = $parts.Count; $i -ge 1; $i--) {
    $part = $parts.Item($i)
    if ($strayNamespaces -contains $part.NamespaceURI -or $strayIds -contains $part.Id) {
        $part.Delete()
    }
}

# Some hosts expose namespace-scoped lookup instead of (or in addition
# to) a fully enumerable collection -- cover that path too.
foreach ($ns in $strayNamespaces) {
    $scoped = $parts.SelectByNamespace($ns)
    for ($i = $scoped.Count; $i -ge 1; $i--) {
        $scoped.Item($i).Delete()
    }
}

# And the equivalent ID-keyed lookup, for hosts that key CustomXML parts
# by their storage GUID rather than by namespace.
foreach ($id in $strayIds) {
    $byId = $parts.SelectByID($id)
    for ($i = $byId.Count; $i -ge 1; $i--) {
        $byId.Item($i).Delete()
    }
}

$d.Save()
